$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Price" (column D) values for the symbol list refresh.
# Leading apostrophe forces the value to be stored as text (matching the
# original inlineStr text cells) rather than being auto-converted to a
# numeric value by Excel's input parsing.
$ws.Range("D2").Value = "'275.66"
$ws.Range("D3").Value = "'21.14"
$ws.Range("D4").Value = "'6.253"
$ws.Range("D5").Value = "'0.06216"
$ws.Range("D6").Value = "'3.554"
$ws.Range("D7").Value = "'1.541"
$ws.Range("D8").Value = "'6.553"
$ws.Range("D10").Value = "'0.1649"
$ws.Range("D11").Value = "'0.08284"
$ws.Range("D12").Value = "'0.03482"
$ws.Range("D13").Value = "'0.03132"
$ws.Range("D14").Value = "'0.09149"
$ws.Range("D16").Value = "'0.001643"
$ws.Range("D17").Value = "'0.04681"
$ws.Range("D18").Value = "'0.006238"
$ws.Range("D19").Value = "'0.006214"
$ws.Range("D20").Value = "'0.001067"
$ws.Range("D21").Value = "'0.0001497"
$ws.Range("D22").Value = "'3.725"
$ws.Range("D24").Value = "'0.01394"
$ws.Range("D25").Value = "'0.3292"
$ws.Range("D26").Value = "'0.1249"
$ws.Range("D28").Value = "'0.0002731"
$ws.Range("D40").Value = "'0.04739"
$ws.Range("D41").Value = "'0.005289"
$ws.Range("D42").Value = "'0.007017"
$ws.Range("D44").Value = "'0.01140"
$ws.Range("D45").Value = "'0.00006158"
$ws.Range("D47").Value = "'0.7215"
$ws.Range("D48").Value = "'0.001392"
$ws.Range("D49").Value = "'0.00001896"
